# Append a new "traded" record to the repeater sheet, and fill in the
# PriceChange / UpDown values for the previous (now-complete) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: fill in the two columns that were left blank (X/Y) ---------
$ws.Range("X10").Value = 0.6499990000000011
$ws.Range("Y10").Value = "Up"

# --- Row 11: brand-new record -------------------------------------------
# Copy formatting from row 10 first so number formats (date, percentages)
# line up with the rest of the table, then overwrite values cell-by-cell.
$ws.Range("A10:W10").Copy($ws.Range("A11:W11"))

$ws.Range("A11").Value = 42654.883287037039
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = "Buy"
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = "Random"
$ws.Range("Q11").Value = 58.438771163779279
$ws.Range("R11").Value = 0.49
$ws.Range("S11").Value = 0.0933
$ws.Range("T11").Value = 0.0249
$ws.Range("U11").Value = 2.34
$ws.Range("V11").Value = "N/A"
$ws.Range("W11").Value = 2

# Row 11 has no entries in the PriceChange / UpDown columns yet, so make
# sure the copy above didn't leave stray values behind.
$ws.Range("X11").ClearContents()
$ws.Range("Y11").ClearContents()

Write-Output "edit applied"
